$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the updated "Price" figures are plain digit-and-dot strings (e.g.
# "587.87") that Excel's value-assignment type-sniffer would otherwise treat
# as numbers, silently dropping significant trailing zeros (e.g. "2.20" ->
# 2.2). The source column is text (coin prices/links/percent strings), so
# force those cells to stay Text the same way a user would in the UI - type
# a leading apostrophe, then reset the cell style back to Normal so no
# left-over custom number format/quote-prefix styling remains on the cell.
function Set-TextCell($Cell, $Text) {
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

$ws.Range('D2').Value = '64.328.02'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '3.306.09'
$ws.Range('E3').Value = '  +5.42%  '
$ws.Range('E4').Value = '  +1.08%  '
Set-TextCell $ws.Range('D5') '587.87'
$ws.Range('E5').Value = '  -0.48%  '
Set-TextCell $ws.Range('D6') '146.84'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.97%  '
$ws.Range('D8').Value = '3.150.44'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  -1.29%  '
Set-TextCell $ws.Range('D11') '5.82'
$ws.Range('E11').Value = '  +2.00%  '
Set-TextCell $ws.Range('D12') '0.458'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('E13').Value = '  -2.75%  '
Set-TextCell $ws.Range('D14') '37.13'
$ws.Range('E14').Value = '  +3.54%  '
$ws.Range('D15').Value = '3.804.42'
$ws.Range('E15').Value = '  +4.31%  '
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('D17').Value = '3.206.05'
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').Value = '64.075.95'
$ws.Range('E18').Value = '  +0.61%  '
Set-TextCell $ws.Range('D19') '7.08'
$ws.Range('E19').Value = '  -1.31%  '
Set-TextCell $ws.Range('D20') '465.05'
$ws.Range('E20').Value = '  -0.31%  '
Set-TextCell $ws.Range('D21') '14.32'
$ws.Range('E21').Value = '  +0.63%  '
Set-TextCell $ws.Range('D22') '0.732'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E23').Value = '  -1.44%  '
Set-TextCell $ws.Range('D24') '12.95'
$ws.Range('E24').Value = '  -2.73%  '
Set-TextCell $ws.Range('D25') '81.18'
$ws.Range('E25').Value = '  -0.87%  '
Set-TextCell $ws.Range('D26') '2.24'
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('E27').Value = '  +0.36%  '
Set-TextCell $ws.Range('D28') '9.29'
$ws.Range('E28').Value = '  +7.72%  '
$ws.Range('E29').Value = '  +0.65%  '
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('E32').Value = '  +3.46%  '
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('E36').Value = '  -0.39%  '
Set-TextCell $ws.Range('D37') '2.33'
$ws.Range('E37').Value = '  -3.23%  '
Set-TextCell $ws.Range('D38') '3.32'
$ws.Range('E38').Value = '  -1.37%  '
Set-TextCell $ws.Range('D39') '6.02'
$ws.Range('E39').Value = '  -1.92%  '
Set-TextCell $ws.Range('D40') '51.44'
$ws.Range('E40').Value = '  +1.30%  '
Set-TextCell $ws.Range('D41') '438.51'
$ws.Range('E41').Value = '  -2.36%  '
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws.Range('D43') '0.288'
$ws.Range('E43').Value = '  +2.95%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range('D44') '0.0372'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.917.41'
$ws.Range('E45').Value = '  -0.35%  '
Set-TextCell $ws.Range('D46') '39.35'
$ws.Range('E46').Value = '  +15.58%  '
Set-TextCell $ws.Range('D47') '0.108'
$ws.Range('E47').Value = '  -3.77%  '
Set-TextCell $ws.Range('D48') '126.57'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E50').Value = '  -1.21%  '
Set-TextCell $ws.Range('D51') '2.20'
$ws.Range('E51').Value = '  +0.43%  '
